# Re-append the Android header + first data row further down the sheet
# (rows 25:26), mirroring existing rows 1:2, and update the
# selection / active-sheet UI state to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Android"
$ws2 = $wb.Worksheets.Item(2)   # "iOS"

# Copy A1:B2 (header "Name"/"Purpose" + "Battle.net - Android.xml"/"Chats")
# down to A25:B26, carrying over formatting (bold header row) as well.
[void]$ws1.Range("A1:B2").Copy($ws1.Range("A25"))

# iOS sheet is no longer the active tab; its lingering selection moves
# from B2 to A2.
[void]$ws2.Range("A2").Select()

# Android becomes the active tab again, scrolled down to the newly
# pasted rows, with the insertion point sitting just below them.
[void]$ws1.Activate()
[void]$excel.Goto($ws1.Range("A25"), $true)
[void]$ws1.Range("A28").Select()
